# Apply updated "dSF" (column F) values for the mahle_tyler save_data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 7
$ws.Range("F8").Value = -3
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -3
